$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '41.530.00'
$ws.Range('E2').Value = '  +0.16%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.482.29'
$ws.Range('E3').Value = '  +0.58%  '
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.39'
$ws.Range('E5').Value = '  +0.47%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '92.36'
$ws.Range('E6').Value = '  -2.69%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.548'
$ws.Range('E7').Value = '  -0.37%  '
$ws.Range('E8').Value = '  -0.22%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.515'
$ws.Range('E9').Value = '  +2.97%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '32.87'
$ws.Range('E10').Value = '  -1.99%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0791'
$ws.Range('E11').Value = '  +1.26%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.110'
$ws.Range('E12').Value = '  +1.83%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.863.79'
$ws.Range('E13').Value = '  +0.53%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.92'
$ws.Range('E14').Value = '  -1.38%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '16.35'
$ws.Range('E15').Value = '  +9.01%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.466.83'
$ws.Range('E16').Value = '  -0.40%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.777'
$ws.Range('E17').Value = '  -1.34%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '41.558.78'
$ws.Range('E18').Value = '  +0.34%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.57'
$ws.Range('E19').Value = '  +3.81%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0946'
$ws.Range('E20').Value = '  +2.39%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.68'
$ws.Range('E21').Value = '  +6.05%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.22'
$ws.Range('E22').Value = '  -0.07%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '237.22'
$ws.Range('E23').Value = '  +0.15%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('B26').Value = 'ImmutableX'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.90'
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.89'
$ws.Range('E27').Value = '  +3.03%  '
$ws.Range('E28').Value = '  +0.18%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.72'
$ws.Range('E29').Value = '  +0.74%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '35.88'
$ws.Range('E30').Value = '  -2.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '157.77'
$ws.Range('E31').Value = '  +3.83%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.47'
$ws.Range('E32').Value = '  -0.38%  '
$ws.Range('E33').Value = '  -0.30%  '
$ws.Range('E34').Value = '  +1.39%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '17.38'
$ws.Range('E35').Value = '  +1.86%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.39'
$ws.Range('E36').Value = '  -10.13%  '
$ws.Range('E37').Value = '  +3.89%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.93'
$ws.Range('E38').Value = '  -4.00%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.84'
$ws.Range('E39').Value = '  -2.28%  '
$ws.Range('E40').Value = '  -0.04%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.08'
$ws.Range('E41').Value = '  -4.03%  '
$ws.Range('E42').Value = '  -0.31%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.972.66'
$ws.Range('E43').Value = '  -0.75%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '19.18'
$ws.Range('E44').Value = '  -3.64%  '
$ws.Range('E45').Value = '  -0.41%  '
$ws.Range('E46').Value = '  -2.63%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.93'
$ws.Range('E47').Value = '  +2.38%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.720.83'
$ws.Range('E48').Value = '  +0.31%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '98.11'
$ws.Range('E49').Value = '  +1.58%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '68.42'
$ws.Range('E50').Value = '  -2.09%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.173'
$ws.Range('E51').Value = '  -2.92%  '
